$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 363
$ws.Range("J3").Value = 8071
$ws.Range("K3").Value = 324
$ws.Range("K4").Value = 67
$ws.Range("K5").Value = 19
$ws.Range("J6").Value = 11049
$ws.Range("K6").Value = 440
$ws.Range("J7").Value = 29216
$ws.Range("K7").Value = 1213

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 21
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 19
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 34
$ws.Range("K8").Value = 81
$ws.Range("K10").Value = 6
$ws.Range("J11").Value = 536
$ws.Range("K11").Value = 35
$ws.Range("K13").Value = 5
$ws.Range("K14").Value = 7
$ws.Range("J16").Value = 111
$ws.Range("K19").Value = 23
$ws.Range("K22").Value = 3
$ws.Range("K25").Value = 4
$ws.Range("K29").Value = 61
$ws.Range("K33").Value = 55
$ws.Range("K35").Value = 6
$ws.Range("K36").Value = 18
$ws.Range("K37").Value = 32
$ws.Range("K41").Value = 14
$ws.Range("K42").Value = 36
$ws.Range("J48").Value = 322
$ws.Range("J49").Value = 176
$ws.Range("K49").Value = 12
$ws.Range("J54").Value = 574
$ws.Range("K54").Value = 19
$ws.Range("K60").Value = 10
$ws.Range("J63").Value = 84
$ws.Range("K65").Value = 34
$ws.Range("J72").Value = 109
$ws.Range("K76").Value = 19
$ws.Range("K77").Value = 12
$ws.Range("K79").Value = 27
$ws.Range("K84").Value = 9
$ws.Range("K85").Value = 60
$ws.Range("K89").Value = 19
$ws.Range("K91").Value = 13
$ws.Range("K94").Value = 11
$ws.Range("K95").Value = 26
$ws.Range("K96").Value = 17
$ws.Range("K98").Value = 5
$ws.Range("J101").Value = 29216
$ws.Range("K101").Value = 1213

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 5
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K2").Value = 1
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 176
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 7
$ws.Range("J4").Value = 47
$ws.Range("J7").Value = 574
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 20
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 62
$ws.Range("J7").Value = 322

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 4
$ws.Range("K4").Value = 2
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 7

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("K4").Value = 1
$ws.Range("K6").Value = 5

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K2").Value = 2
$ws.Range("K3").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 9
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 3
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K5").Value = 1
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 9
$ws.Range("J3").Value = 91
$ws.Range("J7").Value = 536
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 3

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 1
$ws.Range("K7").Value = 3

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 6
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 3

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 111
